# The sheet contains a daily price log for "Ají" (chili pepper) sold at
# "Vega Central Mapocho de Santiago". A new daily record was inserted at
# row 239, pushing the existing rows 239-326 down to 240-327 (dimension
# grows from A1:R326 to A1:R327).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 239 - shifts rows 239..326 down to 240..327
$ws.Rows.Item(239).Insert()

# Populate the newly inserted row 239 with the new record's data
$ws.Range("A239").Value = 9
$ws.Range("B239").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C239").Value = "Metropolitana"
$ws.Range("D239").Value = 44795
$ws.Range("E239").Value = 13
$ws.Range("F239").Value = 100112021
$ws.Range("G239").Value = "Ají"
$ws.Range("H239").Value = "Inferno"
$ws.Range("I239").Value = "Primera"
$ws.Range("J239").Value = 80
$ws.Range("K239").Value = 13000
$ws.Range("L239").Value = 13000
$ws.Range("M239").Value = 13000
$ws.Range("N239").Value = "`$/caja 10 kilos"
$ws.Range("O239").Value = "Región de Arica y Parinacota"
$ws.Range("P239").Value = 1300
$ws.Range("Q239").Value = 10
$ws.Range("R239").Value = "Hortaliza"
